$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 494
$ws1.Range("F3").Value = 5847
$ws1.Range("F5").Value = 78
$ws1.Range("F6").Value = 101
$ws1.Range("F7").Value = 5
$ws1.Range("F8").Value = 58
$ws1.Range("F9").Value = 548
$ws1.Range("F10").Value = 26

# Sheet "全部类型" (sheet4): update "想去人数" (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 494
$ws4.Range("F3").Value = 5847
$ws4.Range("F6").Value = 78
$ws4.Range("F7").Value = 101
$ws4.Range("F8").Value = 5
$ws4.Range("F10").Value = 58
$ws4.Range("F11").Value = 548
$ws4.Range("F12").Value = 26
